$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "MbrNo" in column L, row 1
$ws.Range("L1").Value = "MbrNo"

# Fill member numbers 1..10 down column L for rows 2..11
for ($i = 2; $i -le 11; $i++) {
    $ws.Cells.Item($i, 12).Value = $i - 1
}

# Update the visible selection to N7 (matches the post-edit view state)
$ws.Range("N7").Select()
